$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current (before) data in rows 5-8 (sorted descending by column A: 6, 4, 6, 8)
# Target (after) data: row5 -> old row6 (4, 1, 3, ...), row6 -> old row5 (6, 0, 2, ...) with
# new "merged" style on A6, row7 keeps old row7 B..H values but A7 becomes blank (merged with A6),
# row8 unchanged.

# Step 1: capture old row 5 and row 6 values (B..G only, A handled specially)
$b5 = $ws.Range("B5:G5").Value2
$b6 = $ws.Range("B6:G6").Value2

# Step 2: write swapped values
$ws.Range("A5").Value2 = 4
$ws.Range("B5:G5").Value2 = $b6

$ws.Range("A6").Value2 = 6
$ws.Range("B6:G6").Value2 = $b5

# Step 3: clear A7's value (its number moved into the merged A6 cell), then merge A6:A7
# (merging applies the top-left cell's current style to the whole range, so do this before
# fixing up the individual cell styles below).
$ws.Range("A7").ClearContents()
$ws.Range("A6:A7").Merge()

# Step 4: apply the new style to A6: same as A5/A7/A8 (style 12) plus vertical-top alignment.
# Copy format from A8 (style 12) onto A6 first so the base style matches exactly, then add
# the vertical alignment on top of it.
$ws.Range("A8").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").VerticalAlignment = -4160

# Step 5: give A7 a plain border-only style (same as the other placeholder cells in column H).
$ws.Range("H5").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$excel.CutCopyMode = 0
